# Add the team's season record (Wins / Losses / Ties) as three new
# trailing columns (AD, AE, AF) to the player table on Sheet1.
#
# Every player row shares the same team record, so the same three
# numbers are written down the whole column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Give the new header cells the same look (bold / bordered / centered)
# as the rest of row 1 by copying the formatting from the neighboring
# header cell (AC1) before writing the new header text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$wins = 67
$losses = 95
$ties = 0

$lastRow = 44
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = $wins    # column AD
    $ws.Cells.Item($r, 31).Value = $losses  # column AE
    $ws.Cells.Item($r, 32).Value = $ties    # column AF
}
